# Auto-generated edit script applying cell-value updates per the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '27.163.72'
$ws.Cells.Item(2, 5).Value = '  +0.85%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.901.05'
$ws.Cells.Item(3, 5).Value = '  +1.35%  '
$ws.Cells.Item(4, 5).Value = '  +0.15%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '306.52'
$ws.Cells.Item(5, 5).Value = '  +0.05%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '1.000'
$ws.Cells.Item(6, 5).Value = '  +0.08%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.5235'
$ws.Cells.Item(7, 5).Value = '  +1.53%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3772'
$ws.Cells.Item(8, 5).Value = '  +1.53%  '
$ws.Cells.Item(9, 5).Value = '  +0.80%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '21.14'
$ws.Cells.Item(10, 5).Value = '  +2.25%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.8988'
$ws.Cells.Item(11, 5).Value = '  +0.05%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.08360'
$ws.Cells.Item(12, 5).Value = '  +10.54%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '1.909.88'
$ws.Cells.Item(13, 5).Value = '  +1.77%  '
$ws.Cells.Item(14, 5).Value = '  -0.20%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '5.268'
$ws.Cells.Item(15, 5).Value = '  +0.37%  '
$ws.Cells.Item(16, 5).Value = '  +0.16%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.000008579'
$ws.Cells.Item(17, 5).Value = '  +1.15%  '
$ws.Cells.Item(18, 5).Value = '  +1.75%  '
$ws.Cells.Item(19, 5).Value = '  +0.10%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '27.204.55'
$ws.Cells.Item(20, 5).Value = '  +0.89%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '5.056'
$ws.Cells.Item(21, 5).Value = '  +0.53%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '2.153.12'
$ws.Cells.Item(22, 5).Value = '  +1.79%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '10.58'
$ws.Cells.Item(23, 5).Value = '  +1.82%  '
$ws.Cells.Item(24, 5).Value = '  -0.26%  '
$ws.Cells.Item(25, 5).Value = '  +8.36%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '146.56'
$ws.Cells.Item(26, 5).Value = '  +0.37%  '
$ws.Cells.Item(27, 5).Value = '  -1.47%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '18.11'
$ws.Cells.Item(28, 5).Value = '  +0.45%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '114.65'
$ws.Cells.Item(29, 5).Value = '  +0.11%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '4.923'
$ws.Cells.Item(30, 5).Value = '  +0.54%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '4.777'
$ws.Cells.Item(31, 5).Value = '  +0.54%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.8186'
$ws.Cells.Item(33, 5).Value = '  +8.62%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.05049'
$ws.Cells.Item(34, 5).Value = '  +0.29%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '1.235'
$ws.Cells.Item(35, 5).Value = '  +5.53%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.959'
$ws.Cells.Item(36, 5).Value = '  -1.18%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '3.365'
$ws.Cells.Item(37, 5).Value = '  +2.46%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '2.572'
$ws.Cells.Item(38, 5).Value = '  +3.50%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.5686'
$ws.Cells.Item(39, 5).Value = '  +1.90%  '
$ws.Cells.Item(40, 5).Value = '  -1.05%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.074'
$ws.Cells.Item(41, 5).Value = '  +0.12%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '6.654'
$ws.Cells.Item(42, 5).Value = '  +1.30%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '8.940'
$ws.Cells.Item(43, 5).Value = '  +2.40%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '118.25'
$ws.Cells.Item(44, 5).Value = '  +1.69%  '
$ws.Cells.Item(45, 5).Value = '  +0.49%  '
$ws.Cells.Item(46, 5).Value = '  +0.98%  '
$ws.Cells.Item(47, 2).Value = 'PaxDollar'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.000'
$ws.Cells.Item(47, 5).Value = '  +0.10%  '
$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '10.13'
$ws.Cells.Item(48, 5).Value = '  -0.05%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.607'
$ws.Cells.Item(49, 5).Value = '  +2.85%  '
$ws.Cells.Item(50, 5).Value = '  +0.80%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '63.54'
$ws.Cells.Item(51, 5).Value = '  +0.33%  '
